$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.288150666666667
$ws.Range("H2").Value = 6.864452
$ws.Range("I2").Value = 0.3964219041944151
$ws.Range("J2").Value = 0.3964219041944151
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 20.424575
$ws.Range("N2").Value = 61.273725
$ws.Range("O2").Value = 0.203732656096709
$ws.Range("P2").Value = 0.2037326560967089
$ws.Range("Q2").Value = 46.73450490263333
$ws.Range("R2").Value = 420.6105441237
$ws.Range("S2").Value = 0.08076408747644327
$ws.Range("T2").Value = 0.08076408747644326
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.288150666666667
$ws.Range("H3").Value = 6.864452
$ws.Range("I3").Value = 0.3964219041944151
$ws.Range("J3").Value = 0.3964219041944151
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 57.16769933333334
$ws.Range("N3").Value = 171.503098
$ws.Range("O3").Value = 0.5702408607336045
$ws.Range("P3").Value = 0.5702408607336045
$ws.Range("Q3").Value = 130.8083093413662
$ws.Range("R3").Value = 1177.274784072296
$ws.Range("S3").Value = 0.2260559678614777
$ws.Range("T3").Value = 0.2260559678614777
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.288150666666667
$ws.Range("H4").Value = 6.864452
$ws.Range("I4").Value = 0.3964219041944151
$ws.Range("J4").Value = 0.3964219041944151
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.1182046666666667
$ws.Range("N4").Value = 0.354614
$ws.Range("O4").Value = 0.001179077200040937
$ws.Range("P4").Value = 0.001179077200040937
$ws.Range("Q4").Value = 0.2704700868364445
$ws.Range("R4").Value = 2.434230781528
$ws.Range("S4").Value = 0.0004674120288324474
$ws.Range("T4").Value = 0.0004674120288324474
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.288150666666667
$ws.Range("H5").Value = 6.864452
$ws.Range("I5").Value = 0.3964219041944151
$ws.Range("J5").Value = 0.3964219041944151
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 22.54136766666667
$ws.Range("N5").Value = 67.62410300000001
$ws.Range("O5").Value = 0.2248474059696456
$ws.Range("P5").Value = 0.2248474059696456
$ws.Range("Q5").Value = 51.57804545406179
$ws.Range("R5").Value = 464.202409086556
$ws.Range("S5").Value = 0.08913443682766162
$ws.Range("T5").Value = 0.0891344368276616
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.588894
$ws.Range("H6").Value = 7.766681999999999
$ws.Range("I6").Value = 0.4485256605643812
$ws.Range("J6").Value = 0.4485256605643813
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 20.424575
$ws.Range("N6").Value = 61.273725
$ws.Range("O6").Value = 0.203732656096709
$ws.Range("P6").Value = 0.2037326560967089
$ws.Range("Q6").Value = 52.87705967005
$ws.Range("R6").Value = 475.89353703045
$ws.Range("S6").Value = 0.0913793241543123
$ws.Range("T6").Value = 0.09137932415431228
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.588894
$ws.Range("H7").Value = 7.766681999999999
$ws.Range("I7").Value = 0.4485256605643812
$ws.Range("J7").Value = 0.4485256605643813
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 57.16769933333334
$ws.Range("N7").Value = 171.503098
$ws.Range("O7").Value = 0.5702408607336045
$ws.Range("P7").Value = 0.5702408607336045
$ws.Range("Q7").Value = 148.0011137978707
$ws.Range("R7").Value = 1332.010024180836
$ws.Range("S7").Value = 0.2557676587413413
$ws.Range("T7").Value = 0.2557676587413413
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.588894
$ws.Range("H8").Value = 7.766681999999999
$ws.Range("I8").Value = 0.4485256605643812
$ws.Range("J8").Value = 0.4485256605643813
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.1182046666666667
$ws.Range("N8").Value = 0.354614
$ws.Range("O8").Value = 0.001179077200040937
$ws.Range("P8").Value = 0.001179077200040937
$ws.Range("Q8").Value = 0.3060193523053333
$ws.Range("R8").Value = 2.754174170748
$ws.Range("S8").Value = 0.0005288463800047622
$ws.Range("T8").Value = 0.0005288463800047622
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.588894
$ws.Range("H9").Value = 7.766681999999999
$ws.Range("I9").Value = 0.4485256605643812
$ws.Range("J9").Value = 0.4485256605643813
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 22.54136766666667
$ws.Range("N9").Value = 67.62410300000001
$ws.Range("O9").Value = 0.2248474059696456
$ws.Range("P9").Value = 0.2248474059696456
$ws.Range("Q9").Value = 58.35721150402733
$ws.Range("R9").Value = 525.214903536246
$ws.Range("S9").Value = 0.1008498312887229
$ws.Range("T9").Value = 0.1008498312887229
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.3123523333333333
$ws.Range("H10").Value = 0.9370569999999999
$ws.Range("I10").Value = 0.05411501461132016
$ws.Range("J10").Value = 0.05411501461132018
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 20.424575
$ws.Range("N10").Value = 61.273725
$ws.Range("O10").Value = 0.203732656096709
$ws.Range("P10").Value = 0.2037326560967089
$ws.Range("Q10").Value = 6.379663658591666
$ws.Range("R10").Value = 57.41697292732499
$ws.Range("S10").Value = 0.01102499566147647
$ws.Range("T10").Value = 0.01102499566147647
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.3123523333333333
$ws.Range("H11").Value = 0.9370569999999999
$ws.Range("I11").Value = 0.05411501461132016
$ws.Range("J11").Value = 0.05411501461132018
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 57.16769933333334
$ws.Range("N11").Value = 171.503098
$ws.Range("O11").Value = 0.5702408607336045
$ws.Range("P11").Value = 0.5702408607336045
$ws.Range("Q11").Value = 17.85646427806511
$ws.Range("R11").Value = 160.708178502586
$ws.Range("S11").Value = 0.03085859251057079
$ws.Range("T11").Value = 0.0308585925105708
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.3123523333333333
$ws.Range("H12").Value = 0.9370569999999999
$ws.Range("I12").Value = 0.05411501461132016
$ws.Range("J12").Value = 0.05411501461132018
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.1182046666666667
$ws.Range("N12").Value = 0.354614
$ws.Range("O12").Value = 0.001179077200040937
$ws.Range("P12").Value = 0.001179077200040937
$ws.Range("Q12").Value = 0.03692150344422222
$ws.Range("R12").Value = 0.332293530998
$ws.Range("S12").Value = 0.00006380577990808976
$ws.Range("T12").Value = 0.00006380577990808977
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.3123523333333333
$ws.Range("H13").Value = 0.9370569999999999
$ws.Range("I13").Value = 0.05411501461132016
$ws.Range("J13").Value = 0.05411501461132018
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 22.54136766666667
$ws.Range("N13").Value = 67.62410300000001
$ws.Range("O13").Value = 0.2248474059696456
$ws.Range("P13").Value = 0.2248474059696456
$ws.Range("Q13").Value = 7.040848787207889
$ws.Range("R13").Value = 63.367639084871
$ws.Range("S13").Value = 0.01216762065936481
$ws.Range("T13").Value = 0.01216762065936481
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.5826116666666666
$ws.Range("H14").Value = 1.747835
$ws.Range("I14").Value = 0.1009374206298835
$ws.Range("J14").Value = 0.1009374206298836
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 20.424575
$ws.Range("N14").Value = 61.273725
$ws.Range("O14").Value = 0.203732656096709
$ws.Range("P14").Value = 0.2037326560967089
$ws.Range("Q14").Value = 11.89959568170833
$ws.Range("R14").Value = 107.096361135375
$ws.Range("S14").Value = 0.02056424880447692
$ws.Range("T14").Value = 0.02056424880447692
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.5826116666666666
$ws.Range("H15").Value = 1.747835
$ws.Range("I15").Value = 0.1009374206298835
$ws.Range("J15").Value = 0.1009374206298836
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 57.16769933333334
$ws.Range("N15").Value = 171.503098
$ws.Range("O15").Value = 0.5702408607336045
$ws.Range("P15").Value = 0.5702408607336045
$ws.Range("Q15").Value = 33.30656858809223
$ws.Range("R15").Value = 299.75911729283
$ws.Range("S15").Value = 0.05755864162021468
$ws.Range("T15").Value = 0.05755864162021469
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.5826116666666666
$ws.Range("H16").Value = 1.747835
$ws.Range("I16").Value = 0.1009374206298835
$ws.Range("J16").Value = 0.1009374206298836
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.1182046666666667
$ws.Range("N16").Value = 0.354614
$ws.Range("O16").Value = 0.001179077200040937
$ws.Range("P16").Value = 0.001179077200040937
$ws.Range("Q16").Value = 0.06886741785444445
$ws.Range("R16").Value = 0.61980676069
$ws.Range("S16").Value = 0.0001190130112956374
$ws.Range("T16").Value = 0.0001190130112956374
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.5826116666666666
$ws.Range("H17").Value = 1.747835
$ws.Range("I17").Value = 0.1009374206298835
$ws.Range("J17").Value = 0.1009374206298836
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 22.54136766666667
$ws.Range("N17").Value = 67.62410300000001
$ws.Range("O17").Value = 0.2248474059696456
$ws.Range("P17").Value = 0.2248474059696456
$ws.Range("Q17").Value = 13.13286378522278
$ws.Range("R17").Value = 118.195774067005
$ws.Range("S17").Value = 0.02269551719389631
$ws.Range("T17").Value = 0.02269551719389631
